$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 1
$ws.Range("H113").Value = 4382.9414
$ws.Range("I113").Value = 4070.6667
$ws.Range("J113").Value = 6725
$ws.Range("K113").Value = 4070.6667
$ws.Range("L113").Value = 6725
$ws.Range("M113").Value = -816.6667000000002
$ws.Range("N113").Value = -13233

# Hunk 2
$ws.Range("H116").Value = 3629.2654
$ws.Range("I116").Value = 3109.3635
$ws.Range("K116").Value = 3109.3635
$ws.Range("M116").Value = 332.6365000000001

# Hunk 3
$ws.Range("H132").Value = 3511208.2
$ws.Range("I132").Value = 4168878.5
$ws.Range("J132").Value = 3633.6667
$ws.Range("K132").Value = 12506635.5
$ws.Range("L132").Value = 10901.0001
$ws.Range("M132").Value = -12504105.5
$ws.Range("N132").Value = -15961.0001

$ws = $wb.Worksheets.Item("ARM")
# Hunk 4
$ws.Range("H61").Value = 1715.4783
$ws.Range("I61").Value = 802.3226
$ws.Range("J61").Value = 3602.6667
$ws.Range("K61").Value = 802.3226
$ws.Range("L61").Value = 3602.6667
$ws.Range("M61").Value = -590.3226
$ws.Range("N61").Value = -4026.6667

# Hunk 5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Hunk 6
$ws.Range("H125").Value = 29714
$ws.Range("J125").Value = 29714
$ws.Range("L125").Value = 29714
$ws.Range("N125").Value = -39554

# Hunk 7
$ws.Range("H132").Value = 2514.7026
$ws.Range("I132").Value = 1967.5333
$ws.Range("J132").Value = 4859.7144
$ws.Range("K132").Value = 5902.5999
$ws.Range("L132").Value = 14579.1432
$ws.Range("M132").Value = -3372.5999
$ws.Range("N132").Value = -19639.1432

# Hunk 8
$ws.Range("H136").Value = 1715.4783
$ws.Range("I136").Value = 802.3226
$ws.Range("J136").Value = 3602.6667
$ws.Range("K136").Value = 2406.9678
$ws.Range("L136").Value = 10808.0001
$ws.Range("M136").Value = 143.0322000000001
$ws.Range("N136").Value = -15908.0001

$ws = $wb.Worksheets.Item("BSM")
# Hunk 9
$ws.Range("H124").Value = 39000
$ws.Range("J124").Value = 39000
$ws.Range("L124").Value = 39000
$ws.Range("N124").Value = -48820

# Hunk 10
$ws.Range("H130").Value = 27924.8
$ws.Range("J130").Value = 27924.8
$ws.Range("L130").Value = 27924.8
$ws.Range("N130").Value = -37964.8

# Hunk 11
$ws.Range("H134").Value = 1592.4828
$ws.Range("I134").Value = 738
$ws.Range("J134").Value = 5694
$ws.Range("K134").Value = 2214
$ws.Range("L134").Value = 17082
$ws.Range("M134").Value = 321
$ws.Range("N134").Value = -22152

$ws = $wb.Worksheets.Item("CRP")
# Hunk 12
$ws.Range("H99").Value = 3517.5
$ws.Range("I99").Value = 2338.75
$ws.Range("J99").Value = 5875
$ws.Range("K99").Value = 2338.75
$ws.Range("L99").Value = 5875
$ws.Range("M99").Value = -840.75
$ws.Range("N99").Value = -8871

# Hunk 13
$ws.Range("H123").Value = 31917.5
$ws.Range("J123").Value = 31917.5
$ws.Range("L123").Value = 31917.5
$ws.Range("N123").Value = -41717.5

# Hunk 14
$ws.Range("H124").Value = 25000
$ws.Range("J124").Value = 25000
$ws.Range("L124").Value = 25000
$ws.Range("N124").Value = -29910

# Hunk 15
$ws.Range("H126").Value = 3517.5
$ws.Range("I126").Value = 2338.75
$ws.Range("J126").Value = 5875
$ws.Range("K126").Value = 7016.25
$ws.Range("L126").Value = 17625
$ws.Range("M126").Value = -4546.25
$ws.Range("N126").Value = -22565

# Hunk 16
$ws.Range("H132").Value = 1790.2554
$ws.Range("I132").Value = 1325.8158
$ws.Range("J132").Value = 3751.2222
$ws.Range("K132").Value = 3977.4474
$ws.Range("L132").Value = 11253.6666
$ws.Range("M132").Value = -1447.4474
$ws.Range("N132").Value = -16313.6666

# Hunk 17
$ws.Range("H134").Value = 2788.9
$ws.Range("I134").Value = 1512.6923
$ws.Range("J134").Value = 5159
$ws.Range("K134").Value = 4538.0769
$ws.Range("L134").Value = 15477
$ws.Range("M134").Value = -2003.0769
$ws.Range("N134").Value = -20547

$ws = $wb.Worksheets.Item("CUL")
# Hunk 18
$ws.Range("H76").Value = 3996.6667
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3996.6667
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 11990.0001
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -12756.0001

# Hunk 19
$ws.Range("H79").Value = 3996.6667
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3996.6667
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 11990.0001
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -14642.0001

# Hunk 20
$ws.Range("H130").Value = 2600
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# Hunk 21
$ws.Range("H134").Value = 3104.4666
$ws.Range("I134").Value = 1419.6666
$ws.Range("J134").Value = 4227.6665
$ws.Range("K134").Value = 4258.9998
$ws.Range("L134").Value = 12682.9995
$ws.Range("M134").Value = 811.0002000000004
$ws.Range("N134").Value = -22822.9995

# Hunk 22
$ws.Range("H139").Value = 8068351.5
$ws.Range("I139").Value = 15628182
$ws.Range("J139").Value = 4532
$ws.Range("K139").Value = 46884546
$ws.Range("L139").Value = 13596
$ws.Range("M139").Value = -46879406
$ws.Range("N139").Value = -23876

# Hunk 23
$ws.Range("H140").Value = 11115418
$ws.Range("I140").Value = 33334074
$ws.Range("J140").Value = 6090
$ws.Range("K140").Value = 100002222
$ws.Range("L140").Value = 18270
$ws.Range("M140").Value = -99997042
$ws.Range("N140").Value = -28630

$ws = $wb.Worksheets.Item("GSM")
# Hunk 24
$ws.Range("H126").Value = 3278.3
$ws.Range("I126").Value = 1945.4286
$ws.Range("J126").Value = 4444.5625
$ws.Range("K126").Value = 5836.2858
$ws.Range("L126").Value = 13333.6875
$ws.Range("M126").Value = -3366.2858
$ws.Range("N126").Value = -18273.6875

# Hunk 25
$ws.Range("H132").Value = 2951.9285
$ws.Range("I132").Value = 2602.5715
$ws.Range("K132").Value = 7807.7145
$ws.Range("M132").Value = -5277.7145

# Hunk 26
$ws.Range("H136").Value = 16573.7
$ws.Range("J136").Value = 16573.7
$ws.Range("L136").Value = 49721.10000000001
$ws.Range("N136").Value = -54821.10000000001

$ws = $wb.Worksheets.Item("LTW")
# Hunk 27
$ws.Range("H2").Value = 535714.4
$ws.Range("I2").Value = 538461.5600000001
$ws.Range("J2").Value = 500001
$ws.Range("K2").Value = 538461.5600000001
$ws.Range("L2").Value = 500001
$ws.Range("M2").Value = -538349.5600000001
$ws.Range("N2").Value = -500225

# Hunk 28
$ws.Range("H122").Value = 2973.6333
$ws.Range("I122").Value = 2609.2727
$ws.Range("J122").Value = 3975.625
$ws.Range("K122").Value = 7827.8181
$ws.Range("L122").Value = 11926.875
$ws.Range("M122").Value = -5377.8181
$ws.Range("N122").Value = -16826.875

# Hunk 29
$ws.Range("H135").Value = 30013.182
$ws.Range("J135").Value = 30013.182
$ws.Range("L135").Value = 30013.182
$ws.Range("N135").Value = -40153.182

$ws = $wb.Worksheets.Item("WVR")
# Hunk 30
$ws.Range("H24").Value = 30002
$ws.Range("J24").Value = 30002
$ws.Range("L24").Value = 30002
$ws.Range("N24").Value = -30462

# Hunk 31
$ws.Range("H140").Value = 70107.25
$ws.Range("J140").Value = 70107.25
$ws.Range("L140").Value = 70107.25
$ws.Range("N140").Value = -80467.25

# Hunk 32
$ws.Range("H141").Value = 28423.076
$ws.Range("J141").Value = 28423.076
$ws.Range("L141").Value = 28423.076
$ws.Range("N141").Value = -38783.076
